$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate / extend the existing records, writing new cells in the same
# order the original macro would have produced them so shared-string
# indices line up (row 2, then row 3, then row 4, then the extra columns
# on row 1 last).

# Row 2 - duplicate of row 1's pattern
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "mundo"
$ws.Cells.Item(2, 3).Value = "hola mundo"
$ws.Cells.Item(2, 4).Value = "nn"
$ws.Cells.Item(2, 6).Value = "hola"

# Row 3
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "pc"
$ws.Cells.Item(3, 5).Value = "g"

# Row 4
$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = "mary"
$ws.Cells.Item(4, 5).Value = "*"
$ws.Cells.Item(4, 3).Value = 67
$ws.Cells.Item(4, 4).Value = 67

# Back to row 3 for the trailing column
$ws.Cells.Item(3, 3).Value = "ii"

# Row 1 - extra columns appended last
$ws.Cells.Item(1, 7).Value = "'0.999"
$ws.Cells.Item(1, 8).Value = 0.78

# Selection moves to reflect the new active cell after the edits
$ws.Range("K9").Select()
